$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Logistic Regression): fill in the "Count Vectorizer" configuration
# and its associated results.
$ws.Range("B2").Value = "Count Vectorizer"
$ws.Range("C2").Value = 82.46
$ws.Range("D2").Value = 84.42
$ws.Range("E2").Value = 92.86
$ws.Range("F2").Value = 88.32

# Rows 3-6 keep their "NULL" configuration (no count-vectorization results yet).
$ws.Range("B3").Value = "NULL"
$ws.Range("B4").Value = "NULL"
$ws.Range("B5").Value = "NULL"
$ws.Range("B6").Value = "NULL"

# Update the active selection to match the new state of the sheet.
$ws.Range("F3").Select()
